$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move last_name column (header + values) from column B to column I,
# and move personal_email values into column D, clearing out the stale
# duplicate email columns/values that caused NullPointerException on import.

$ws.Range("I1").Value2 = $ws.Range("B1").Value2
$ws.Range("I2").Value2 = $ws.Range("B2").Value2
$ws.Range("I3").Value2 = $ws.Range("B3").Value2

$ws.Range("B1:B3").ClearContents()

$ws.Range("D2").Value2 = "john@personal.com"
$ws.Range("D3").Value2 = "mary@personal.com"

$ws.Range("G2").ClearContents()

$ws.Columns.Item(9).ColumnWidth = 17.42
$ws.Columns.Item(2).ColumnWidth = $ws.StandardWidth
$ws.Columns.Item(4).ColumnWidth = 18.59

$ws.Range("C7").Select()
